$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.99999999550202245,
    0.99283520683621196,
    0.96440315485856054,
    0.95611271779620932,
    0.9482797127603122,
    0.92948079187027233,
    0.92710507076234405,
    0.92437089066789158,
    0.92361497048004682,
    0.92381089084298063,
    0.92395469048301837,
    0.92117468633073518,
    0.90988698371531918,
    0.905719458043581,
    0.90312787562476338,
    0.90062131934552792,
    0.89691331828487675,
    0.89580439366589248,
    0.98933922890635295,
    0.98222201557100564,
    0.98082351065542106,
    0.97955900265734641,
    0.9481672462432682,
    0.92453875419235543,
    0.90771492526571174,
    0.92717153535628383,
    0.92420173684317208,
    0.91143684017547977,
    0.90263429410829854,
    0.89928806939768502,
    0.89986374028544702,
    0.89818440793113397,
    0.89766439260319442
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
